$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 82, shifting rows 82:138 down to 83:139
$ws.Rows("82:82").Insert()

# Fill in values for the new row 82
$ws.Cells.Item(82, 1).Value = 6
$ws.Cells.Item(82, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(82, 3).Value = "Metropolitana"
$ws.Cells.Item(82, 4).Value = 45126
$ws.Cells.Item(82, 5).Value = 13
$ws.Cells.Item(82, 6).Value = 100114007
$ws.Cells.Item(82, 7).Value = "Jengibre"
$ws.Cells.Item(82, 8).Value = "Sin especificar"
$ws.Cells.Item(82, 9).Value = "Primera"
$ws.Cells.Item(82, 10).Value = 430
$ws.Cells.Item(82, 11).Value = 16000
$ws.Cells.Item(82, 12).Value = 17000
$ws.Cells.Item(82, 13).Value = 16814
$ws.Cells.Item(82, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(82, 15).Value = "Perú"
$ws.Cells.Item(82, 16).Value = 1293
$ws.Cells.Item(82, 17).Value = 13
$ws.Cells.Item(82, 18).Value = "Hortaliza"
